# Update "Training Dashboard" sheet with a new progress-as-of date (04-Nov-2025):
#  - Column H (PERIOD TO EXPIRE) decreases by 1 for every data row
#  - Column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 43; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    $jCell = $ws.Cells.Item($row, 10)  # column J - untouched cell used as a format donor

    $hVal = $hCell.Value2
    if ($hVal -ne $null) {
        $hCell.Value2 = $hVal - 1
    }

    $iVal = $iCell.Value2
    if ($iVal -eq "03-Nov-2025") {
        # Writing a dash-separated date-looking string straight into Value2
        # makes Excel auto-convert it into a date serial (and re-style the
        # cell). Force the cell to Text first so the literal string
        # "04-Nov-2025" is kept, then restore the original cell style (copied
        # from an untouched neighbour cell that still carries style s="3").
        $iCell.NumberFormat = "@"
        $iCell.Value2 = "04-Nov-2025"
        $jCell.Copy()
        $iCell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

$excel.CutCopyMode = 0
